{"js": "// Apply benchmark-stat corrections to the single-column results table.\n// Each row's cell is addressed by (row, col) index so we don't have to\n// worry about duplicate text values elsewhere in the table (e.g. \"103\").\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// row index (0-based) -> new text for column 0\nconst updates = [\n  [0, \"0M\"],\n  [1, \"0M\"],\n  [2, \"0M\"],\n  [3, \"412\"],\n  [5, \"0.00076\"],\n  [6, \"0.00025\"],\n  [7, \"0.00006\"],\n  [8, \"0.00040\"],\n  [9, \"0.00048\"],\n  [10, \"0.00055\"],\n  [11, \"0.10443\"],\n  // The last three rows previously held a whole tab-separated summary\n  // line crammed into one run; they are replaced with the single\n  // percentage value that used to live in row 0/1/2.\n  [43, \"99.88\"],\n  [44, \"0.1\"],\n  [45, \"90\"],\n];\n\nfor (const [rowIndex, newText] of updates) {\n  const cell = table.getCell(rowIndex, 0);\n  const range = cell.body.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply benchmark-stat corrections to the single-column results table.\n# Cells are addressed by their (row, column) position (1-based, as COM\n# expects) so duplicate text values elsewhere in the table (e.g. \"103\")\n# are not accidentally touched.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"412\"\n    6  = \"0.00076\"\n    7  = \"0.00025\"\n    8  = \"0.00006\"\n    9  = \"0.00040\"\n    10 = \"0.00048\"\n    11 = \"0.00055\"\n    12 = \"0.10443\"\n    # The last three rows previously held a whole tab-separated summary\n    # line crammed into one run; they are replaced with the single\n    # percentage value that used to live in row 1/2/3.\n    44 = \"99.88\"\n    45 = \"0.1\"\n    46 = \"90\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $cell = $t.Cell($rowIndex, 1)\n    $cell.Range.Text = $updates[$rowIndex]\n}\n"}
